$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update C2 text (longer note, keeps same shared-string slot) ---
$ws.Range("C2").Value = "working on getting min-devkit environment setup. Big compiler issues with cmake. Had to do this: had to do this: https://stackoverflow.com/questions/17980759/xcode-select-active-developer-directory-error/17980786#17980786. "

# --- Update B2 value 2 -> 4 ---
$ws.Range("B2").Value = 4

# --- Add new row 3 data ---
# A3 must stay literal text "2020.02.22" (not get auto-converted to a date serial).
# Route it through a formula + paste-values so Excel's live-typing date-sniffer never sees it.
$ws.Range("A3").Formula = "=""2020.02.22"""
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B3").Value = 1.5
$ws.Range("C3").Value = "Went through two buffer Max examples. Have yet to go through corresponding code throroughly. Read through GuideToAudio.md. Read through c74_min_operator_vector.h, which contains critical information about the audio_bundle object."

# --- Wrap text for the notes column (creates cellXfs index 1) ---
$ws.Range("C1:C3").WrapText = $true

# --- Column C width ---
$ws.Columns.Item(3).ColumnWidth = 77.66666666666667

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 54
$ws.Rows.Item(3).RowHeight = 45

# --- Update selection to A4 ---
[void]$ws.Range("A4").Select()

Write-Output "done"
